# Update "descriptives" sheet with revised descriptive statistics
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("descriptives")
$ws1.Range("C2").Value = 845
$ws1.Range("D2").Value = 0.007710466534301554
$ws1.Range("E2").Value = 0.05059402236771077
$ws1.Range("F2").Value = 97.7210678207355
$ws1.Range("G2").Value = 84.79796297379484
$ws1.Range("H2").Value = 12.92310484694066
$ws1.Range("C3").Value = 438
$ws1.Range("D3").Value = 0.002266906743304885
$ws1.Range("E3").Value = 0.05323324433729858
$ws1.Range("F3").Value = 94.98187742734611
$ws1.Range("G3").Value = 91.10233738575046
$ws1.Range("H3").Value = 3.879540041595643

# Update "coefficients" sheet with revised model coefficients
$ws2 = $wb.Worksheets.Item("coefficients")
$ws2.Range("D2").Value = 0.1088660576832843
$ws2.Range("E2").Value = 0.02694838947763055
$ws2.Range("F2").Value = 4.055872400564315
$ws2.Range("G2").Value = 0.0003777071052438656
$ws2.Range("H2").Value = 0.05397136768315828
$ws2.Range("I2").Value = 0.1631046937958562
$ws2.Range("J2").Value = 27.19094419875941
$ws2.Range("D3").Value = 0.04570203060089384
$ws2.Range("E3").Value = 0.04905655216241486
$ws2.Range("F3").Value = 0.9322687271233721
$ws2.Range("G3").Value = 0.4086086668163802
$ws2.Range("H3").Value = -0.09543082549685647
$ws2.Range("I3").Value = 0.1850337295565198
$ws2.Range("J3").Value = 3.653246753684463
$ws2.Range("E4").Value = [double]"9.675137993899536e-17"
$ws2.Range("F4").Value = 1870047346286094
$ws2.Range("G4").Value = [double]"3.40429761648498e-16"
$ws2.Range("H4").Value = 0.1789808979992103
$ws2.Range("D5").Value = 0.1366960143109622
$ws2.Range("E5").Value = [double]"3.207913979655756e-16"
$ws2.Range("F5").Value = 428805494825438.2
$ws2.Range("G5").Value = [double]"1.484635295139447e-15"
$ws2.Range("H5").Value = 0.1366960143109582
$ws2.Range("I5").Value = 0.1366960143109662
$ws2.Range("D6").Value = 0.0776500003536626
$ws2.Range("E6").Value = 0.02221418991573963
$ws2.Range("F6").Value = 3.502564434267127
$ws2.Range("G6").Value = 0.001942747956487021
$ws2.Range("H6").Value = 0.03180972527748716
$ws2.Range("I6").Value = 0.1231642938999653
$ws2.Range("J6").Value = 22.70889626145427
$ws2.Range("D7").Value = 0.4572828288958798
$ws2.Range("E7").Value = [double]"2.166474313636758e-16"
$ws2.Range("F7").Value = 2279603726948210
$ws2.Range("G7").Value = [double]"2.7926773624811e-16"
$ws2.Range("H7").Value = 0.4572828288958776
$ws2.Range("I7").Value = 0.457282828895882
$ws2.Range("D8").Value = 0.172930865397675
$ws2.Range("E8").Value = 0.04636330590028399
$ws2.Range("F8").Value = 3.767770799168988
$ws2.Range("G8").Value = 0.003370664097575282
$ws2.Range("H8").Value = 0.07192822190898303
$ws2.Range("I8").Value = 0.2704228398527846
$ws2.Range("J8").Value = 10.50581634719451
$ws2.Range("D9").Value = 0.09243101045783822
$ws2.Range("E9").Value = 0.0302921192009756
$ws2.Range("F9").Value = 3.060056489278443
$ws2.Range("G9").Value = 0.01001196620783359
$ws2.Range("H9").Value = 0.02660735758199492
$ws2.Range("I9").Value = 0.1574565946474072
$ws2.Range("J9").Value = 11.86865584060162
$ws2.Range("D10").Value = 0.05229526606655354
$ws2.Range("E10").Value = 0.03157866084237695
$ws2.Range("F10").Value = 1.657543902197481
$ws2.Range("G10").Value = 0.1688235084642944
$ws2.Range("H10").Value = -0.03345555782361513
$ws2.Range("I10").Value = 0.1372817791184651
$ws2.Range("J10").Value = 4.232258305185446
$ws2.Range("D11").Value = 0.05126955589916281
$ws2.Range("E11").Value = [double]"2.361288975329441e-17"
$ws2.Range("F11").Value = 2173158357308966
$ws2.Range("G11").Value = [double]"2.929467934200207e-16"
$ws2.Range("H11").Value = 0.05126955589916251
$ws2.Range("I11").Value = 0.05126955589916311
$ws2.Range("D12").Value = -0.005095168936260564
$ws2.Range("E12").Value = 0.031285224964531
$ws2.Range("F12").Value = -0.1628632376525016
$ws2.Range("G12").Value = 0.8732287039183231
$ws2.Range("H12").Value = -0.07282729001473054
$ws2.Range("I12").Value = 0.06268373525099136
$ws2.Range("J12").Value = 12.50386954169758
$ws2.Range("D13").Value = 0.05693802332938057
$ws2.Range("E13").Value = 0.0387315676428654
$ws2.Range("F13").Value = 1.471659336773607
$ws2.Range("G13").Value = 0.1747822883382157
$ws2.Range("H13").Value = -0.0304418403151962
$ws2.Range("I13").Value = 0.1434542099698444
$ws2.Range("J13").Value = 9.113442660591378

# Update "pairwise" sheet with revised pairwise comparison statistics
$ws3 = $wb.Worksheets.Item("pairwise")
$ws3.Range("C2").Value = 1.135681993343508
$ws3.Range("E2").Value = 4.907697647035894
$ws3.Range("F2").Value = 0.3084886948449619
$ws3.Range("G2").Value = 0.3966283219435224
$ws3.Range("C3").Value = 2.658059853332448
$ws3.Range("E3").Value = 27.1909441987594
$ws3.Range("F3").Value = 0.01300454816233488
$ws3.Range("G3").Value = 0.03385258680827456
$ws3.Range("C4").Value = 1.048592719867346
$ws3.Range("E4").Value = 27.1909441987594
$ws3.Range("F4").Value = 0.3035971223524434
$ws3.Range("G4").Value = 0.3966283219435224
$ws3.Range("C5").Value = 0.9017465241997387
$ws3.Range("E5").Value = 47.8956291057825
$ws3.Range("F5").Value = 0.37170397851715
$ws3.Range("G5").Value = 0.44604477422058
$ws3.Range("C6").Value = 14.27065105941262
$ws3.Range("E6").Value = 27.1909441987594
$ws3.Range("F6").Value = [double]"3.778049700936994e-14"
$ws3.Range("G6").Value = [double]"3.400244730843294e-13"
$ws3.Range("C7").Value = 1.219311689927433
$ws3.Range("E7").Value = 19.39394332933332
$ws3.Range("F7").Value = 0.2373484712518597
$ws3.Range("G7").Value = 0.3560227068777896
$ws3.Range("C8").Value = 2.75591671047974
$ws3.Range("E8").Value = 3.653246753684462
$ws3.Range("F8").Value = 0.05647248700513487
$ws3.Range("G8").Value = 0.1270630957615535
$ws3.Range("C9").Value = 1.871783088115357
$ws3.Range("E9").Value = 3.653246753684462
$ws3.Range("F9").Value = 0.1413046943486819
$ws3.Range("G9").Value = 0.2312258634796613
$ws3.Range("C10").Value = 0.5955742728264293
$ws3.Range("E10").Value = 5.247762442424044
$ws3.Range("F10").Value = 0.5762175311730369
$ws3.Range("G10").Value = 0.6101126800655684
$ws3.Range("C11").Value = 9.135097816193861
$ws3.Range("E11").Value = 3.653246753684462
$ws3.Range("F11").Value = 0.001205866535826753
$ws3.Range("G11").Value = 0.004341119528976311
$ws3.Range("C12").Value = 1.910437377367318
$ws3.Range("E12").Value = 6.803424116407768
$ws3.Range("F12").Value = 0.09890122967779887
$ws3.Range("G12").Value = 0.178022213420038
$ws3.Range("C14").Value = 4.642214285583472
$ws3.Range("E14").Value = 22.70889626145426
$ws3.Range("F14").Value = 0.0001168870880583285
$ws3.Range("G14").Value = 0.0005259918962624783
$ws3.Range("C16").Value = 0.1346614754744741
$ws3.Range("E16").Value = 10.50581634719452
$ws3.Range("F16").Value = 0.8954246335837814
$ws3.Range("G16").Value = 0.8954246335837814
$ws3.Range("C17").Value = 2.689743935389337
$ws3.Range("E17").Value = 22.70889626145426
$ws3.Range("F17").Value = 0.01316489486988455
$ws3.Range("G17").Value = 0.03385258680827456
$ws3.Range("C19").Value = 0.800831503765332
$ws3.Range("E19").Value = 10.50581634719452
$ws3.Range("F19").Value = 0.4409640780905094
$ws3.Range("G19").Value = 0.4960845878518231
$ws3.Range("C20").Value = 18.7296346170458
$ws3.Range("E20").Value = 22.70889626145426
$ws3.Range("F20").Value = [double]"2.650709969264459e-15"
$ws3.Range("G20").Value = [double]"4.771277944676026e-14"
$ws3.Range("C21").Value = 1.884438382418104
$ws3.Range("E21").Value = 21.24446019660779
$ws3.Range("F21").Value = 0.0732624463167395
$ws3.Range("G21").Value = 0.146524892633479
$ws3.Range("C22").Value = 6.884409462363371
$ws3.Range("E22").Value = 10.50581634719452
$ws3.Range("F22").Value = [double]"3.33583406293437e-05"
$ws3.Range("G22").Value = 0.0002001500437760622
$ws3.Range("C23").Value = 0.9221611737124695
$ws3.Range("E23").Value = 8.517108786496292
$ws3.Range("F23").Value = 0.3818242952866268
$ws3.Range("G23").Value = 0.6363738254777113
$ws3.Range("C24").Value = 1.36606643308144
$ws3.Range("E24").Value = 11.8686558406016
$ws3.Range("F24").Value = 0.1972441186778474
$ws3.Range("G24").Value = 0.4655763641206606
$ws3.Range("C25").Value = 2.245616594784637
$ws3.Range("E25").Value = 24.10139663986222
$ws3.Range("F25").Value = 0.03416480890029303
$ws3.Range("G25").Value = 0.3416480890029303
$ws3.Range("C26").Value = 0.7259608547951426
$ws3.Range("E26").Value = 20.25085496711703
$ws3.Range("F26").Value = 0.4761700227506654
$ws3.Range("G26").Value = 0.6802428896438077
$ws3.Range("C27").Value = 0.03256844919392969
$ws3.Range("E27").Value = 4.232258305185445
$ws3.Range("F27").Value = 0.9754973396540498
$ws3.Range("G27").Value = 0.9754973396540498
$ws3.Range("C28").Value = 1.292141545926189
$ws3.Range("E28").Value = 7.90754487392106
$ws3.Range("F28").Value = 0.2327881820603303
$ws3.Range("G28").Value = 0.4655763641206606
$ws3.Range("C29").Value = 0.09318260022456785
$ws3.Range("E29").Value = 8.778216188743157
$ws3.Range("F29").Value = 0.9278505087117879
$ws3.Range("G29").Value = 0.9754973396540498
$ws3.Range("C30").Value = 1.803079945883427
$ws3.Range("E30").Value = 12.50386954169757
$ws3.Range("F30").Value = 0.09551270829270293
$ws3.Range("G30").Value = 0.4655763641206606
$ws3.Range("C31").Value = 0.1467827095028937
$ws3.Range("E31").Value = 9.11344266059138
$ws3.Range("F31").Value = 0.8864993540887041
$ws3.Range("G31").Value = 0.9754973396540498
$ws3.Range("C32").Value = 1.247172033613706
$ws3.Range("E32").Value = 19.8658827403374
$ws3.Range("F32").Value = 0.2268380661295306
$ws3.Range("G32").Value = 0.4655763641206606

# Add new "nr_studies" sheet: number of effect sizes and studies per moderator level
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "nr_studies"

$ws4.Range("A1").Value = "outcome"
$ws4.Range("B1").Value = "moderator_context"
$ws4.Range("C1").Value = "n_effect_sizes"
$ws4.Range("D1").Value = "k_studies"
$headerRange = $ws4.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

$ws4.Range("A2").Value = "NS"
$ws4.Range("B2").Value = "Education"
$ws4.Range("C2").Value = 417
$ws4.Range("D2").Value = 33
$ws4.Range("A3").Value = "NS"
$ws4.Range("B3").Value = "Health care"
$ws4.Range("C3").Value = 10
$ws4.Range("D3").Value = 1
$ws4.Range("A4").Value = "NS"
$ws4.Range("B4").Value = "Parenting"
$ws4.Range("C4").Value = 240
$ws4.Range("D4").Value = 28
$ws4.Range("A5").Value = "NS"
$ws4.Range("B5").Value = "Sport"
$ws4.Range("C5").Value = 132
$ws4.Range("D5").Value = 13
$ws4.Range("A6").Value = "NS"
$ws4.Range("B6").Value = "Exercise"
$ws4.Range("C6").Value = 37
$ws4.Range("D6").Value = 5
$ws4.Range("A7").Value = "NS"
$ws4.Range("B7").Value = "Healthcare"
$ws4.Range("C7").Value = 6
$ws4.Range("D7").Value = 1
$ws4.Range("A8").Value = "NS"
$ws4.Range("B8").Value = "Partner"
$ws4.Range("C8").Value = 3
$ws4.Range("D8").Value = 1
$ws4.Range("A9").Value = "NT"
$ws4.Range("B9").Value = "Education"
$ws4.Range("C9").Value = 140
$ws4.Range("D9").Value = 15
$ws4.Range("A10").Value = "NT"
$ws4.Range("B10").Value = "Exercise"
$ws4.Range("C10").Value = 44
$ws4.Range("D10").Value = 6
$ws4.Range("A11").Value = "NT"
$ws4.Range("B11").Value = "Health care"
$ws4.Range("C11").Value = 10
$ws4.Range("D11").Value = 1
$ws4.Range("A12").Value = "NT"
$ws4.Range("B12").Value = "Parenting"
$ws4.Range("C12").Value = 129
$ws4.Range("D12").Value = 17
$ws4.Range("A13").Value = "NT"
$ws4.Range("B13").Value = "Sport"
$ws4.Range("C13").Value = 115
$ws4.Range("D13").Value = 12
